$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range('D2').Value = '61.797.37'
$ws.Range('E2').Value = '  -0.98%  '

$ws.Range('D3').Value = '2.397.08'
$ws.Range('E3').Value = '  -1.24%  '

$ws.Range('E4').Value = '  -0.04%  '

$ws.Range('D5').Value = '559.18'
$ws.Range('E5').Value = '  +0.41%  '

$ws.Range('D6').Value = '141.41'
$ws.Range('E6').Value = '  -1.66%  '

$ws.Range('E8').Value = '  -0.65%  '

$ws.Range('E9').Value = '  -1.78%  '

$ws.Range('E10').Value = '  -1.88%  '

$ws.Range('D11').Value = '5.20'
$ws.Range('E11').Value = '  -3.80%  '

$ws.Range('E12').Value = '  -1.48%  '

$ws.Range('D13').Value = '25.23'
$ws.Range('E13').Value = '  -4.25%  '

$ws.Range('E14').Value = '  -2.74%  '

$ws.Range('D15').Value = '2.829.83'
$ws.Range('E15').Value = '  -1.20%  '

$ws.Range('D16').Value = '61.699.30'
$ws.Range('E16').Value = '  -0.78%  '

$ws.Range('D17').Value = '2.397.97'
$ws.Range('E17').Value = '  -1.13%  '

$ws.Range('D18').Value = '11.12'
$ws.Range('E18').Value = '  -0.08%  '

$ws.Range('D19').Value = '319.82'
$ws.Range('E19').Value = '  -1.51%  '

$ws.Range('D20').Value = '6.77'
$ws.Range('E20').Value = '  +0.30%  '

$ws.Range('E21').Value = '  -1.97%  '

$ws.Range('D22').Value = '1.00'
$ws.Range('E22').Value = '  -0.36%  '

$ws.Range('D23').Value = '65.35'
$ws.Range('E23').Value = '  +0.46%  '

$ws.Range('D24').Value = '1.70'
$ws.Range('E24').Value = '  -5.13%  '

$ws.Range('D25').Value = '8.63'
$ws.Range('E25').Value = '  -5.14%  '

$ws.Range('D26').Value = '559.28'
$ws.Range('E26').Value = '  -1.75%  '

$ws.Range('B27').Value = 'Binance-PegBSC-USD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D27').Value = '1.02'
$ws.Range('E27').Value = '  +2.14%  '

$ws.Range('B28').Value = 'WrappedeETH'
$ws.Range('C28').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D28').Value = '2.519.40'
$ws.Range('E28').Value = '  -1.00%  '

$ws.Range('D29').Value = '0.0₃0920'
$ws.Range('E29').Value = '  -2.70%  '

$ws.Range('D30').Value = '8.09'
$ws.Range('E30').Value = '  -3.88%  '

$ws.Range('E31').Value = '  -5.99%  '

$ws.Range('E32').Value = '  -1.91%  '

$ws.Range('E33').Value = '  -0.51%  '

$ws.Range('E34').Value = '  -5.29%  '

$ws.Range('D36').Value = '4.72'
$ws.Range('E36').Value = '  -2.59%  '

$ws.Range('D37').Value = '151.90'
$ws.Range('E37').Value = '  +1.31%  '

$ws.Range('B38').Value = 'PolygonEcosystemToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D38').Value = '0.377'
$ws.Range('E38').Value = '  -2.16%  '

$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D39').Value = '5.38'
$ws.Range('E39').Value = '  -6.19%  '

$ws.Range('D40').Value = '18.39'
$ws.Range('E40').Value = '  -2.30%  '

$ws.Range('E41').Value = '  -6.45%  '

$ws.Range('E42').Value = '  -0.06%  '

$ws.Range('D43').Value = '146.72'
$ws.Range('E43').Value = '  -3.19%  '

$ws.Range('E44').Value = '  -5.80%  '

$ws.Range('D45').Value = '3.57'
$ws.Range('E45').Value = '  -1.96%  '

$ws.Range('D46').Value = '0.0525'
$ws.Range('E46').Value = '  -3.60%  '

$ws.Range('D47').Value = '19.67'
$ws.Range('E47').Value = '  -3.77%  '

$ws.Range('D48').Value = '0.584'
$ws.Range('E48').Value = '  -1.26%  '

$ws.Range('D49').Value = '0.0913'
$ws.Range('E49').Value = '  -0.38%  '

$ws.Range('E50').Value = '  -2.46%  '

$ws.Range('D51').Value = '11.53'
$ws.Range('E51').Value = '  +0.30%  '
